$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.713.85"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "2.050.03"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.81"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.32"
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E9").Value = "  -1.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0839"
$ws.Range("E10").Value = "  +3.27%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "2.352.65"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.42"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.49"
$ws.Range("E15").Value = "  +5.92%  "
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").Value = "2.056.73"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").Value = "37.726.32"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.95"
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.43"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "0.0₃0831"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "222.81"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("E25").Value = "  +4.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.12"
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.31"
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.129"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.81"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.29"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("E32").Value = "  +8.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.38"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0604"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("E37").Value = "  +3.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.44"
$ws.Range("E38").Value = "  +6.46%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.37"
$ws.Range("E40").Value = "  +9.47%  "
$ws.Range("D41").Value = "1.525.45"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.20"
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("E45").Value = "  +1.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0890"
$ws.Range("E46").Value = "  -2.99%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Value = "2.241.69"
$ws.Range("E51").Value = "  +1.05%  "
